# JIRAMETRICS-25 fix to PPT for module info
#
# The "Requires only Python + two easily "pip install"-able modules" bullet
# on the "Quick and Dirty Jira Automation" slide is reworded to call out the
# `jira` module by name:
#
#   Requires only Python + jira module (easily "pip install"-able)
#
# with "jira" split into its own run (it is flagged by the spell checker in
# the authored deck, hence its own <a:r>).

$p = $ppt.ActivePresentation

$openQuote  = [char]8220   # “
$closeQuote = [char]8221   # ”

# NOTE: TextRange.Text round-trips curly quotes as plain ASCII quotes when
# read back through this COM shim (the OOXML on disk keeps the real
# typographic quotes either way), so the match text below uses a plain
# quote even though the file itself has “ / ”.
$oldText = 'Requires only Python + two easily "pip install"-able modules'

$prefix = "Requires only Python + "
$word   = "jira"
$suffix = " module (easily " + $openQuote + "pip install" + $closeQuote + "-able)"

# Locate the slide/shape/paragraph that still has the old wording (robust to
# shape re-ordering) rather than hard-coding indices.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }

        $tf = $shape.TextFrame
        $tr = $tf.TextRange
        if ($tr.Length -le 0) { continue }

        $paraCount = $tr.Paragraphs().Count

        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi, 1)

            if ($para.Text.TrimEnd("`r") -eq $oldText) {
                # Replace whole paragraph text first...
                $para.Text = $prefix + $word + $suffix

                # ...then re-stamp each piece through Characters() so the
                # paragraph ends up as three distinct runs (matching the
                # "jira" run being its own <a:r>), and keep the original
                # run's formatting (sz="2000", dirty="0") on all three.
                $r1 = $para.Characters(1, $prefix.Length)
                $r1.Text = $prefix
                $r1.Font.Size = 20

                $r2 = $para.Characters($prefix.Length + 1, $word.Length)
                $r2.Text = $word
                $r2.Font.Size = 20

                $r3 = $para.Characters($prefix.Length + $word.Length + 1, $suffix.Length)
                $r3.Text = $suffix
                $r3.Font.Size = 20
            }
        }
    }
}
